# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Pepino ensalada" (Macroferia Regional
# de Talca) right before the current row 298, shifting the existing rows
# 298-318 down to 299-319 (dimension grows from A1:R318 to A1:R319).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 298.. down by one, creating a blank row 298 to fill in.
$ws.Rows.Item(298).Insert()

$ws.Cells.Item(298, 1).Value  = 5
$ws.Cells.Item(298, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(298, 3).Value  = "Maule"
$ws.Cells.Item(298, 4).Value  = 44585
$ws.Cells.Item(298, 5).Value  = 7
$ws.Cells.Item(298, 6).Value  = 100112043
$ws.Cells.Item(298, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(298, 8).Value  = "Sin especificar"
$ws.Cells.Item(298, 9).Value  = "Primera"
$ws.Cells.Item(298, 10).Value = 400
$ws.Cells.Item(298, 11).Value = 9000
$ws.Cells.Item(298, 12).Value = 9000
$ws.Cells.Item(298, 13).Value = 9000
$ws.Cells.Item(298, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(298, 15).Value = "Región del Maule"
$ws.Cells.Item(298, 16).Value = 112
$ws.Cells.Item(298, 17).Value = 80
$ws.Cells.Item(298, 18).Value = "Hortaliza"

Write-Output "inserted row 298"
